# Update the "Config" table (Tabelle1) with a new set of DFU package
# assignments (switching away from the old BLE DFU zip files to the
# master/client/server DFU zip files) and shrink the table by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Config")

# New Number / Firmware values for rows 2-8 (Dev ID recalculates via the
# existing VLOOKUP formula already present in column C).
$ws.Range("A2").Value = 45
$ws.Range("B2").Value = "master_dfu_package.zip"

$ws.Range("A3").Value = 10
$ws.Range("B3").Value = "client_dfu_package.zip"

$ws.Range("A4").Value = 11
$ws.Range("B4").Value = "server_dfu_package.zip"

$ws.Range("A5").Value = 13
$ws.Range("B5").Value = "client_dfu_package.zip"

$ws.Range("A6").Value = 14
$ws.Range("B6").Value = "server_dfu_package.zip"

$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "client_dfu_package.zip"

$ws.Range("A8").Value = 15
$ws.Range("B8").Value = "server_dfu_package.zip"

# Row 9 is dropped from the table data entirely.
$ws.Range("A9:C9").ClearContents()

# Group ID (D) and Node Id (E) columns are no longer populated.
$ws.Range("D2:D9").ClearContents()
$ws.Range("E2:E9").ClearContents()

# Shrink the table (Tabelle1) from A1:E9 to A1:E8.
$lo = $ws.ListObjects.Item("Tabelle1")
$lo.Resize($ws.Range("A1:E8"))

# Update the saved cursor/selection position on the Config sheet.
$ws.Range("G16").Select() | Out-Null
